$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 82.98768099999999
$ws.Range("H2").Value = 248.963043
$ws.Range("I2").Value = 0.4489504115427952
$ws.Range("J2").Value = 0.4489504115427952
$ws.Range("M2").Value = 0.4102596666666667
$ws.Range("N2").Value = 1.230779
$ws.Range("O2").Value = 0.003499619873322347
$ws.Range("P2").Value = 0.003499619873322347
$ws.Range("Q2").Value = 34.04649834449967
$ws.Range("R2").Value = 306.418485100497
$ws.Range("S2").Value = 0.001571155782371413
$ws.Range("T2").Value = 0.001571155782371413
$ws.Range("G3").Value = 82.98768099999999
$ws.Range("H3").Value = 248.963043
$ws.Range("I3").Value = 0.4489504115427952
$ws.Range("J3").Value = 0.4489504115427952
$ws.Range("O3").Value = 0.8692174743460166
$ws.Range("P3").Value = 0.8692174743460165
$ws.Range("Q3").Value = 8456.293075406813
$ws.Range("R3").Value = 76106.63767866133
$ws.Range("S3").Value = 0.3902355428278332
$ws.Range("T3").Value = 0.3902355428278332
$ws.Range("G4").Value = 82.98768099999999
$ws.Range("H4").Value = 248.963043
$ws.Range("I4").Value = 0.4489504115427952
$ws.Range("J4").Value = 0.4489504115427952
$ws.Range("N4").Value = 44.764041
$ws.Range("O4").Value = 0.1272829057806611
$ws.Range("P4").Value = 0.1272829057806611
$ws.Range("Q4").Value = 1238.287984926307
$ws.Range("R4").Value = 11144.59186433676
$ws.Range("S4").Value = 0.05714371293259066
$ws.Range("T4").Value = 0.05714371293259066
$ws.Range("G5").Value = 63.14058933333333
$ws.Range("I5").Value = 0.3415807409566563
$ws.Range("J5").Value = 0.3415807409566563
$ws.Range("M5").Value = 0.4102596666666667
$ws.Range("N5").Value = 1.230779
$ws.Range("O5").Value = 0.003499619873322347
$ws.Range("P5").Value = 0.003499619873322347
$ws.Range("Q5").Value = 25.90403713303022
$ws.Range("R5").Value = 233.136334197272
$ws.Range("S5").Value = 0.001195402749396087
$ws.Range("T5").Value = 0.001195402749396087
$ws.Range("G6").Value = 63.14058933333333
$ws.Range("I6").Value = 0.3415807409566563
$ws.Range("J6").Value = 0.3415807409566563
$ws.Range("O6").Value = 0.8692174743460166
$ws.Range("P6").Value = 0.8692174743460165
$ws.Range("Q6").Value = 6433.910695209955
$ws.Range("R6").Value = 57905.19625688959
$ws.Range("S6").Value = 0.2969079489395857
$ws.Range("T6").Value = 0.2969079489395857
$ws.Range("G7").Value = 63.14058933333333
$ws.Range("I7").Value = 0.3415807409566563
$ws.Range("J7").Value = 0.3415807409566563
$ws.Range("N7").Value = 44.764041
$ws.Range("O7").Value = 0.1272829057806611
$ws.Range("P7").Value = 0.1272829057806611
$ws.Range("Q7").Value = 942.1426432271652
$ws.Range("R7").Value = 8479.283789044488
$ws.Range("S7").Value = 0.04347738926767451
$ws.Range("T7").Value = 0.04347738926767451
$ws.Range("I8").Value = 0.2094688475005485
$ws.Range("J8").Value = 0.2094688475005485
$ws.Range("M8").Value = 0.4102596666666667
$ws.Range("N8").Value = 1.230779
$ws.Range("O8").Value = 0.003499619873322347
$ws.Range("P8").Value = 0.003499619873322347
$ws.Range("Q8").Value = 15.88523049827267
$ws.Range("R8").Value = 142.967074484454
$ws.Range("S8").Value = 0.0007330613415548476
$ws.Range("T8").Value = 0.0007330613415548475
$ws.Range("I9").Value = 0.2094688475005485
$ws.Range("J9").Value = 0.2094688475005485
$ws.Range("O9").Value = 0.8692174743460166
$ws.Range("P9").Value = 0.8692174743460165
$ws.Range("R9").Value = 35509.4221361937
$ws.Range("S9").Value = 0.1820739825785977
$ws.Range("T9").Value = 0.1820739825785976
$ws.Range("I10").Value = 0.2094688475005485
$ws.Range("J10").Value = 0.2094688475005485
$ws.Range("N10").Value = 44.764041
$ws.Range("O10").Value = 0.1272829057806611
$ws.Range("P10").Value = 0.1272829057806611
$ws.Range("S10").Value = 0.02666180358039599
$ws.Range("T10").Value = 0.02666180358039599
